$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.791.41'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.597.00'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.86'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.477'
$ws.Range('E7').Value = '  -5.42%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.246'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0609'
$ws.Range('E9').Value = '  -2.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.82'
$ws.Range('E10').Value = '  -3.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0784'
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('D12').Value = '1.819.52'
$ws.Range('E12').Value = '  -2.31%  '
$ws.Range('D13').Value = '1.613.24'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.05'
$ws.Range('E14').Value = '  -3.09%  '
$ws.Range('E15').Value = '  -4.48%  '
$ws.Range('D16').Value = '25.786.67'
$ws.Range('E16').Value = '  -0.79%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '60.25'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('E18').Value = '  -3.94%  '
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '188.57'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('E21').Value = '  -1.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.29'
$ws.Range('E22').Value = '  -3.31%  '
$ws.Range('E23').Value = '  -3.48%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.38'
$ws.Range('E25').Value = '  -1.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.126'
$ws.Range('E26').Value = '  -5.31%  '
$ws.Range('E27').Value = '  -2.78%  '
$ws.Range('E28').Value = '  -4.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '14.88'
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.18'
$ws.Range('E30').Value = '  -4.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0465'
$ws.Range('E31').Value = '  -3.89%  '
$ws.Range('E32').Value = '  -2.73%  '
$ws.Range('E33').Value = '  -5.41%  '
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('E35').Value = '  -2.10%  '
$ws.Range('D36').Value = '1.097.00'
$ws.Range('E36').Value = '  -3.47%  '
$ws.Range('E37').Value = '  -2.73%  '
$ws.Range('B38').Value = 'PaxDollar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.01'
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.794'
$ws.Range('E39').Value = '  -8.13%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0151'
$ws.Range('E40').Value = '  -2.77%  '
$ws.Range('B41').Value = 'ImmutableX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.493'
$ws.Range('E41').Value = '  -5.67%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '95.49'
$ws.Range('E42').Value = '  -3.10%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.732.22'
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.08'
$ws.Range('E44').Value = '  -2.96%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.739'
$ws.Range('E45').Value = '  -4.93%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  -1.45%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '53.03'
$ws.Range('E47').Value = '  -4.02%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0509'
$ws.Range('E48').Value = '  -3.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.43'
$ws.Range('E49').Value = '  -4.00%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.410'
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  -0.21%  '
